$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 values get rounded from 3 decimal places to 2 decimal places
# (custom accuracy). Cells not listed here (M5, N5, Q5, AB5) keep their
# original value since rounding to 2 decimals does not change them.
$row5 = @{
    "B5"  = 11.91
    "C5"  = 8.789999999999999
    "D5"  = 0.95
    "E5"  = 26.12
    "F5"  = 20.98
    "G5"  = 9.31
    "H5"  = 38.23
    "I5"  = 14.54
    "J5"  = 6.37
    "K5"  = 9.31
    "L5"  = 10.46
    "O5"  = 9.4
    "P5"  = 13.3
    "R5"  = 0.72
    "S5"  = 0.59
    "T5"  = 135.74
    "U5"  = 26.34
    "V5"  = 8.68
    "W5"  = 17.54
    "X5"  = 9.17
    "Y5"  = 1.53
    "Z5"  = 18.4
    "AA5" = 7.66
    "AC5" = 8.08
    "AD5" = 10.97
    "AE5" = 0.54
    "AF5" = 34.93
    "AG5" = 4.81
    "AH5" = 10.85
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove the last data row (row 6); this also shrinks the used range /
# sheet dimension from A1:AH6 down to A1:AH5.
$ws.Rows.Item(6).Delete()
